$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (D1, E1 new, with same style as C1)
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column C updated values (rows 2-11)
$cValues = @(
    -5.041423880421429,
    -1.155700669620174,
    -0.07506705564893197,
    -0.4238602485323116,
    0.01520034001876744,
    0.1043279679824023,
    0.1352696695087812,
    0.02787891322180851,
    0.02331057633078736,
    0.006819328375059187
)

# Column D new values (rows 2-11)
$dValues = @(
    -4.710269109749491,
    -1.158930277919864,
    0.01963012362906111,
    -0.1608612870209168,
    -0.04412375307091955,
    0.06603273543210547,
    0.03973136016071439,
    0.02431123640232687,
    0.01663020176682151,
    0.00929739404265134
)

# Column E new values (rows 2-11)
$eValues = @(
    -4.344855546796872,
    -1.140201035762955,
    0.09302013488457662,
    0.04823066005975032,
    -0.09856047570660458,
    0.02905822296111348,
    -0.05081791680323088,
    0.02094463452186093,
    0.01184093569365686,
    0.01118403691153653
)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
